$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'26.007.27"
$ws.Range('E2').Value = '  -0.55%  '
$ws.Range('D3').Value = "'1.648.66"
$ws.Range('E3').Value = '  -0.32%  '
$ws.Range('D4').Value = "'1.001"
$ws.Range('E4').Value = '  -0.40%  '
$ws.Range('D5').Value = "'218.23"
$ws.Range('E5').Value = '  -0.08%  '
$ws.Range('D6').Value = "'0.5272"
$ws.Range('E6').Value = '  +1.30%  '
$ws.Range('D7').Value = "'1.002"
$ws.Range('E7').Value = '  -0.39%  '
$ws.Range('D8').Value = "'0.2618"
$ws.Range('E8').Value = '  -1.85%  '
$ws.Range('D9').Value = "'0.06288"
$ws.Range('E9').Value = '  -0.59%  '
$ws.Range('D10').Value = "'20.30"
$ws.Range('E10').Value = '  -3.70%  '
$ws.Range('D11').Value = "'0.07746"
$ws.Range('E11').Value = '  +0.23%  '
$ws.Range('D12').Value = "'4.468"
$ws.Range('E12').Value = '  +0.46%  '
$ws.Range('D13').Value = "'1.606.13"
$ws.Range('E13').Value = '  -2.87%  '
$ws.Range('D14').Value = "'0.5447"
$ws.Range('E14').Value = '  -0.25%  '
$ws.Range('D15').Value = "'0.0₅8092"
$ws.Range('E15').Value = '  -1.63%  '
$ws.Range('D16').Value = "'64.91"
$ws.Range('E16').Value = '  +0.11%  '
$ws.Range('D17').Value = "'26.014.26"
$ws.Range('E17').Value = '  -0.65%  '
$ws.Range('D18').Value = "'1.002"
$ws.Range('E18').Value = '  -0.35%  '
$ws.Range('D19').Value = "'4.554"
$ws.Range('E19').Value = '  -2.38%  '
$ws.Range('D20').Value = "'192.25"
$ws.Range('E20').Value = '  -0.15%  '
$ws.Range('E21').Value = '  -1.21%  '
$ws.Range('D22').Value = "'5.979"
$ws.Range('E22').Value = '  -1.92%  '
$ws.Range('D23').Value = "'1.003"
$ws.Range('E23').Value = '  -0.51%  '
$ws.Range('D24').Value = "'139.36"
$ws.Range('E24').Value = '  +1.56%  '
$ws.Range('D25').Value = "'0.1243"
$ws.Range('E25').Value = '  +0.53%  '
$ws.Range('E26').Value = '  +0.31%  '
$ws.Range('D27').Value = "'16.20"
$ws.Range('E27').Value = '  +0.62%  '
$ws.Range('D28').Value = "'1.418"
$ws.Range('E28').Value = '  +0.46%  '
$ws.Range('D29').Value = "'0.05939"
$ws.Range('E29').Value = '  -1.59%  '
$ws.Range('E30').Value = '  -0.70%  '
$ws.Range('D31').Value = "'3.490"
$ws.Range('E31').Value = '  -2.31%  '
$ws.Range('D32').Value = "'3.242"
$ws.Range('E32').Value = '  -2.84%  '
$ws.Range('D33').Value = "'1.536"
$ws.Range('E33').Value = '  -6.95%  '
$ws.Range('D34').Value = "'2.410"
$ws.Range('E34').Value = '  -0.06%  '
$ws.Range('D35').Value = "'0.9421"
$ws.Range('E35').Value = '  -3.89%  '
$ws.Range('D36').Value = "'2.749"
$ws.Range('E36').Value = '  -0.86%  '
$ws.Range('D37').Value = "'0.5657"
$ws.Range('E37').Value = '  -4.33%  '
$ws.Range('D38').Value = "'0.01606"
$ws.Range('E38').Value = '  +0.87%  '
$ws.Range('D39').Value = "'5.857"
$ws.Range('E39').Value = '  -1.61%  '
$ws.Range('D40').Value = "'0.8481"
$ws.Range('E40').Value = '  -1.82%  '
$ws.Range('D41').Value = "'1.001"
$ws.Range('E41').Value = '  -0.36%  '
$ws.Range('D42').Value = "'100.81"
$ws.Range('E42').Value = '  +1.02%  '
$ws.Range('D43').Value = "'1.006.41"
$ws.Range('E43').Value = '  -3.04%  '
$ws.Range('D44').Value = "'1.786.60"
$ws.Range('E44').Value = '  -0.42%  '
$ws.Range('D45').Value = "'56.72"
$ws.Range('E45').Value = '  -0.74%  '
$ws.Range('E46').Value = '  -3.09%  '
$ws.Range('D47').Value = "'1.006"
$ws.Range('E47').Value = '  +0.12%  '
$ws.Range('B48').Value = 'RenderToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D48').Value = "'1.482"
$ws.Range('E48').Value = '  +1.07%  '
$ws.Range('B49').Value = 'Mantle'
$ws.Range('C49').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D49').Value = "'0.4283"
$ws.Range('E49').Value = '  +1.26%  '
$ws.Range('D50').Value = "'0.05148"
$ws.Range('E50').Value = '  -0.56%  '
$ws.Range('D51').Value = "'7.827"
$ws.Range('E51').Value = '  -3.35%  '
